$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: PanelID (D1) and PanelPosition (E1), styled like the
# existing header row (A1:C1).
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("D1").Value = "PanelID"
$ws.Range("E1").Value = "PanelPosition"

# New data: which dialogue panel / side each line should appear on.
$ws.Range("E3").Value = "下"
$ws.Range("E4").Value = "下"

# Column C got narrower to make room for the new columns; column E gets an
# explicit width too. (Values pre-compensated for this engine's column-width
# round-trip quantization so the saved width lands as close as possible to
# the target 70.4727272727273 / 17.9272727272727 character widths.)
$ws.Columns("C").ColumnWidth = 69.71428571428571
$ws.Columns("E").ColumnWidth = 17.142857142857142

$ws.Range("D7").Select()
